$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.847.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.526.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0806"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.923.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.525.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.850"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.691.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0936"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "284.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.143"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0776"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.003.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.48%  "
